# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new sheet "PO Forecast" right after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(44969.99999999999, 77, -184.3636069839142, 342.6444210233852),
    @(44976.99999999999, 78, -173.0471121862446, 356.6460805123269),
    @(44983.99999999999, 78, -182.310516282795, 344.4774907833351),
    @(44990.99999999999, 79, -182.9682518788344, 347.3745869635646),
    @(44997.99999999999, 80, -186.9675077483778, 362.9642351760211),
    @(45004.99999999999, 81, -206.4444393148276, 327.7685831817484),
    @(45011.99999999999, 82, -167.0089361083433, 339.4482485793996),
    @(45025.99999999999, 84, -177.8423623169037, 337.2376318384395),
    @(45039.99999999999, 85, -182.1002863593034, 354.4294344016369),
    @(45088.99999999999, 91, -139.4538854204271, 356.7804401487633),
    @(45095.99999999999, 92, -166.7982865604031, 353.0263186302702),
    @(45102.99999999999, 93, -159.7260470160595, 346.7000253631815),
    @(45116.99999999999, 95, -175.8968938248443, 358.8891285774415),
    @(45123.99999999999, 95, -168.2934922354964, 362.1327570208349),
    @(45130.99999999999, 96, -175.2687731924192, 363.2057705576863),
    @(45137.99999999999, 97, -154.3878645965564, 362.6623170482814),
    @(45158.99999999999, 100, -170.4791888923167, 363.5253210577522),
    @(45172.99999999999, 101, -149.6482605123202, 366.0479057115564),
    @(45179.99999999999, 102, -175.8867911444887, 329.2661859849793),
    @(45186.99999999999, 103, -144.1059093322575, 365.5295207208952),
    @(45193.99999999999, 104, -164.3705690489207, 383.2211848098984),
    @(45200.99999999999, 105, -140.602223473324, 370.6704785721655),
    @(45207.99999999999, 106, -140.3545032482949, 359.7298082236087),
    @(45214.99999999999, 106, -157.5876007119241, 387.3030316825586),
    @(45221.99999999999, 107, -165.5078330820832, 382.0698332817145),
    @(45242.99999999999, 110, -151.2168265068043, 376.6094559213005),
    @(45249.99999999999, 111, -149.1874063139441, 373.0731827367729),
    @(45263.99999999999, 112, -146.6787102904629, 345.0707999541996),
    @(45312.99999999999, 118, -146.5118834395173, 362.7409827482956),
    @(45326.99999999999, 120, -116.9059121060404, 383.1177306912695),
    @(45333.99999999999, 121, -165.528853383691, 383.690658332517),
    @(45354.99999999999, 123, -161.7693485605801, 382.2398819821188),
    @(45361.99999999999, 124, -131.5257189584611, 386.9837508589515),
    @(45424.99999999999, 132, -131.4119837281218, 396.4589823911274),
    @(45431.99999999999, 133, -149.0062380137816, 409.2858632783738),
    @(45445.99999999999, 134, -137.306520032451, 417.70006657069),
    @(45459.99999999999, 136, -124.7325210624856, 393.3959782657089),
    @(45466.99999999999, 137, -128.7258627091911, 402.5852463672857),
    @(45473.99999999999, 138, -124.5146810350561, 406.7295798633868),
    @(45487.99999999999, 139, -131.7121303255668, 412.0889324777673),
    @(45494.99999999999, 140, -112.6921477425219, 412.2872216726652),
    @(45571.99999999999, 150, -106.7545104796367, 406.5140582135217),
    @(45592.99999999999, 152, -111.6718871750757, 405.0745600127524),
    @(45599.99999999999, 153, -92.95064839909371, 414.5980312491081),
    @(45613.99999999999, 155, -96.69499504534147, 419.3072427197307),
    @(45620.99999999999, 156, -120.0788071569056, 410.3118103375381),
    @(45627.99999999999, 156, -99.781323007083, 423.5786033373604),
    @(45634.99999999999, 157, -113.9572646392545, 434.6301487583696),
    @(45641.99999999999, 158, -120.6954258350592, 421.779502295145),
    @(45648.99999999999, 159, -115.0965167893619, 434.9905151738202),
    @(45655.99999999999, 160, -100.2445879073675, 419.9495101354744),
    @(45662.99999999999, 161, -82.09208038493128, 437.6492403723318),
    @(45669.99999999999, 161, -103.2505762690654, 427.6214964330471)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}
$lastRow = $r - 1

# --- Match formatting of the other sheets: bold/centered header, date-formatted first column ---
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A" + $lastRow).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
$wsForecast.Range("A1").Select() | Out-Null
